# Update countries & provincias Spain
# - Reorders two pairs of countries (Republica de Africa Central / Costa Rica,
#   and Fiyi / Dominica) to reflect their new ranking position.
# - Refreshes the covid-19 counters for a handful of countries to the
#   "14:18" data refresh (was "13:01").
# - Updates the "datos actualizados" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-rank countries: swap the country names between row pairs ------
# Costa Rica (row 96) <-> Republica de Africa Central (row 97)
$ws.Range("A96").Value = "Republica de Africa Central"
$ws.Range("A97").Value = "Costa Rica"

# Dominica (row 205) <-> Fiyi (row 206)
$ws.Range("A205").Value = "Fiyi"
$ws.Range("A206").Value = "Dominica"

# --- 2. Refresh the "datos actualizados" timestamp ------------------------
$ws.Range("A1").Value = "Datos actualizados a 2 de Julio de 2020 a las 14:18"

# --- 3. Refresh the numeric counters (Casos totales, Nuevos casos,
#        Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
#        columns B..H, keyed by row number. -------------------------------
$changes = @{
    4   = @{ B=2781085; C=1132; D=1165159; E=1485113; F=0; G=15; H=130813 }
    17  = @{ B=196372;  C=48;   D=180300;  E=7011;    F=0; G=0;  H=9061   }
    72  = @{ B=8996;    C=215;  D=5847;    E=3122;    F=0; G=1;  H=27     }
    96  = @{ B=3788;    C=43;   D=810;     E=2931;    F=0; G=0;  H=47     }
    97  = @{ B=3753;    C=0;    D=1516;    E=2220;    F=0; G=0;  H=17     }
    102 = @{ B=2912;    C=81;   D=2155;    E=647;     F=0; G=2;  H=110    }
    107 = @{ D=1954;    E=418;  F=0;       G=1;       H=10 }
    111 = @{ B=2059;    C=5;    D=1827;    E=221;     F=0; G=0;  H=11     }
    115 = @{ B=1850;    C=3;    D=1828;    E=12;      F=0; G=0;  H=10     }
    117 = @{ B=1796;    C=8;    D=1242;    E=519;     F=0; G=1;  H=35     }
    127 = @{ B=1243;    C=9;    D=1120;    E=116;     F=0; G=0;  H=7      }
    136 = @{ B=967;     C=5;    D=846;     E=68;      F=0; G=0;  H=53     }
}

foreach ($r in $changes.Keys) {
    $rowvals = $changes[$r]
    foreach ($col in $rowvals.Keys) {
        $addr = "" + $col + $r
        $ws.Range($addr).Value = $rowvals[$col]
    }
}
